$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.072214007652827
$ws.Range("D2").Value = 1.075219186007867
$ws.Range("E2").Value = 1.085045630645789
$ws.Range("F2").Value = 1.090363757551108
$ws.Range("I2").Value = 1.062857048207258
$ws.Range("J2").Value = 1.07713468491442
$ws.Range("K2").Value = 1.077906333905298
$ws.Range("L2").Value = 1.087707076622973
$ws.Range("M2").Value = 1.093011501941762
$ws.Range("N2").Value = 1.029212821801204

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.073334487314781
$ws.Range("D3").Value = 1.076113494583663
$ws.Range("E3").Value = 1.086112594503774
$ws.Range("F3").Value = 1.091413070642086
$ws.Range("I3").Value = 1.063249379241664
$ws.Range("J3").Value = 1.077912574307286
$ws.Range("K3").Value = 1.078617726725623
$ws.Range("L3").Value = 1.088592561733469
$ws.Range("M3").Value = 1.093880368806591
$ws.Range("N3").Value = 1.029484172785184

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.074059662519852
$ws.Range("D4").Value = 1.076692260955553
$ws.Range("E4").Value = 1.086803440952647
$ws.Range("F4").Value = 1.092092467506469
$ws.Range("I4").Value = 1.063502067416755
$ws.Range("J4").Value = 1.078415461713581
$ws.Range("K4").Value = 1.079077495111189
$ws.Range("L4").Value = 1.089165372127153
$ws.Range("M4").Value = 1.09444239325274
$ws.Range("N4").Value = 1.029659352790288

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.07436456252383
$ws.Range("D5").Value = 1.076935595398119
$ws.Range("E5").Value = 1.087093980359969
$ws.Range("F5").Value = 1.092378186652599
$ws.Range("I5").Value = 1.063608015949515
$ws.Range("J5").Value = 1.078626765725625
$ws.Range("K5").Value = 1.079270649945143
$ws.Range("L5").Value = 1.089406143864626
$ws.Range("M5").Value = 1.094678622674009
$ws.Range("N5").Value = 1.029732902049263

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.074415758720844
$ws.Range("D6").Value = 1.076976453550117
$ws.Range("E6").Value = 1.087142769510967
$ws.Range("F6").Value = 1.092426166056549
$ws.Range("I6").Value = 1.063625788663191
$ws.Range("J6").Value = 1.078662238158496
$ws.Range("K6").Value = 1.079303073761955
$ws.Range("L6").Value = 1.089446568285009
$ws.Range("M6").Value = 1.094718283954744
$ws.Range("N6").Value = 1.029745245636144

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.074063736465223
$ws.Range("D7").Value = 1.076695512319814
$ws.Range("E7").Value = 1.086807322731839
$ws.Range("F7").Value = 1.092096284904118
$ws.Range("I7").Value = 1.063503484212607
$ws.Range("J7").Value = 1.078418285599719
$ws.Range("K7").Value = 1.079080076573447
$ws.Range("L7").Value = 1.089168589479228
$ws.Range("M7").Value = 1.094445549941275
$ws.Range("N7").Value = 1.029660335938079

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.07259264822544
$ws.Range("D8").Value = 1.075521403228715
$ws.Range("E8").Value = 1.085406122978635
$ws.Range("F8").Value = 1.090718290556159
$ws.Range("I8").Value = 1.062989881875375
$ws.Range("J8").Value = 1.077397671350956
$ws.Range("K8").Value = 1.078146866364585
$ws.Range("L8").Value = 1.08800636310946
$ws.Range("M8").Value = 1.093305178985891
$ws.Range("N8").Value = 1.029304609347154

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.070001530667964
$ws.Range("D9").Value = 1.073453163741743
$ws.Range("E9").Value = 1.082940474197007
$ws.Range("F9").Value = 1.088293318439924
$ws.Range("I9").Value = 1.062075839362646
$ws.Range("J9").Value = 1.075595700619024
$ws.Range("K9").Value = 1.076498217374842
$ws.Range("L9").Value = 1.08595716379074
$ws.Range("M9").Value = 1.091294242200214
$ws.Range("N9").Value = 1.02867469831846

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.068274845956974
$ws.Range("D10").Value = 1.072074812452669
$ws.Range("E10").Value = 1.081299023315009
$ws.Range("F10").Value = 1.086678845400419
$ws.Range("I10").Value = 1.061460417250571
$ws.Range("J10").Value = 1.074392015875048
$ws.Range("K10").Value = 1.075396285039141
$ws.Range("L10").Value = 1.084590215645054
$ws.Range("M10").Value = 1.089952637364935
$ws.Range("N10").Value = 1.028252696205633

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.067527337018589
$ws.Range("D11").Value = 1.071478085232975
$ws.Range("E11").Value = 1.080588804915362
$ws.Range("F11").Value = 1.085980276411212
$ws.Range("I11").Value = 1.061192494080942
$ws.Range("J11").Value = 1.073870242720911
$ws.Range("K11").Value = 1.074918463423322
$ws.Range("L11").Value = 1.08399811709139
$ws.Range("M11").Value = 1.089371474193344
$ws.Range("N11").Value = 1.028069476729762

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.067249701713208
$ws.Range("D12").Value = 1.071256450415289
$ws.Range("E12").Value = 1.080325079201635
$ws.Range("F12").Value = 1.085720873049759
$ws.Range("I12").Value = 1.06109275877062
$ws.Range("J12").Value = 1.073676347074558
$ws.Range("K12").Value = 1.074740877228052
$ws.Range("L12").Value = 1.083778154893199
$ws.Range("M12").Value = 1.089155568158769
$ws.Range("N12").Value = 1.028001347320307

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.067309254409539
$ws.Range("D13").Value = 1.071303991100811
$ws.Range("E13").Value = 1.080381645599327
$ws.Range("F13").Value = 1.08577651250228
$ws.Range("I13").Value = 1.06111416214282
$ws.Range("J13").Value = 1.073717942238798
$ws.Range("K13").Value = 1.074778974686834
$ws.Range("L13").Value = 1.083825338912281
$ws.Range("M13").Value = 1.089201882375204
$ws.Range("N13").Value = 1.028015964636936

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.067504387137078
$ws.Range("D14").Value = 1.071459764498389
$ws.Range("E14").Value = 1.080567003613947
$ws.Range("F14").Value = 1.085958832493395
$ws.Range("I14").Value = 1.061184254350997
$ws.Range("J14").Value = 1.073854216999475
$ws.Range("K14").Value = 1.074903786174163
$ws.Range("L14").Value = 1.083979935569633
$ws.Range("M14").Value = 1.089353628071815
$ws.Range("N14").Value = 1.028063846631417

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.06762461784829
$ws.Range("D15").Value = 1.071555743782445
$ws.Range("E15").Value = 1.080681219534467
$ws.Range("F15").Value = 1.086071175959232
$ws.Range("I15").Value = 1.061227411751756
$ws.Range("J15").Value = 1.073938168991157
$ws.Range("K15").Value = 1.074980673133645
$ws.Range("L15").Value = 1.084075183632985
$ws.Range("M15").Value = 1.089447118807167
$ws.Range("N15").Value = 1.028093338565353

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.068324458903074
$ws.Range("D16").Value = 1.072114417549549
$ws.Range("E16").Value = 1.081346169544661
$ws.Range("F16").Value = 1.086725217825725
$ws.Range("I16").Value = 1.061478168046984
$ws.Range("J16").Value = 1.07442663221436
$ws.Range("K16").Value = 1.075427982194914
$ws.Range("L16").Value = 1.084629507021262
$ws.Range("M16").Value = 1.089991202190331
$ws.Range("N16").Value = 1.028264845563635

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.068763492171191
$ws.Range("D17").Value = 1.072464887766442
$ws.Range("E17").Value = 1.081763419723397
$ws.Range("F17").Value = 1.087135617054283
$ws.Range("I17").Value = 1.061635074838094
$ws.Range("J17").Value = 1.074732879567634
$ws.Range("K17").Value = 1.075708385810831
$ws.Range("L17").Value = 1.084977165460676
$ws.Range("M17").Value = 1.090332427089881
$ws.Range("N17").Value = 1.028372296277123

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.069019588042815
$ws.Range("D18").Value = 1.072669321532178
$ws.Range("E18").Value = 1.082006846873346
$ws.Range("F18").Value = 1.087375045051542
$ws.Range("I18").Value = 1.061726456858241
$ws.Range("J18").Value = 1.074911453516167
$ws.Range("K18").Value = 1.075871875124133
$ws.Range("L18").Value = 1.085179929369966
$ws.Range("M18").Value = 1.090531434712056
$ws.Range("N18").Value = 1.028434923201019

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.069106912717932
$ws.Range("D19").Value = 1.072739029929749
$ws.Range("E19").Value = 1.082089858097692
$ws.Range("F19").Value = 1.087456692138442
$ws.Range("I19").Value = 1.061757592184465
$ws.Range("J19").Value = 1.074972333305445
$ws.Range("K19").Value = 1.075927609671851
$ws.Range("L19").Value = 1.085249063343772
$ws.Range("M19").Value = 1.090599287257833
$ws.Range("N19").Value = 1.028456269346273

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.06871638649037
$ws.Range("D20").Value = 1.072427284551206
$ws.Range("E20").Value = 1.081718647332179
$ws.Range("F20").Value = 1.087091580026252
$ws.Range("I20").Value = 1.061618254608796
$ws.Range("J20").Value = 1.074700027810654
$ws.Range("K20").Value = 1.075678307919359
$ws.Range("L20").Value = 1.084939867006131
$ws.Range("M20").Value = 1.090295819276879
$ws.Range("N20").Value = 1.028360772720369

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.067446924789307
$ws.Range("D21").Value = 1.071413892666875
$ws.Range("E21").Value = 1.0805124180636
$ws.Range("F21").Value = 1.085905141695968
$ws.Range("I21").Value = 1.061163619926475
$ws.Range("J21").Value = 1.073814089851475
$ws.Range("K21").Value = 1.074867035118225
$ws.Range("L21").Value = 1.083934411549225
$ws.Range("M21").Value = 1.089308943745896
$ws.Range("N21").Value = 1.02804974860712

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.066648895347342
$ws.Range("D22").Value = 1.070776826942595
$ws.Range("E22").Value = 1.079754482517352
$ws.Range("F22").Value = 1.085159622010472
$ws.Range("I22").Value = 1.06087651912338
$ws.Range("J22").Value = 1.073256568573227
$ws.Range("K22").Value = 1.074356365649916
$ws.Range("L22").Value = 1.083302065610784
$ws.Range("M22").Value = 1.08868824665094
$ws.Range("N22").Value = 1.027853769955795

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.067071933647297
$ws.Range("D23").Value = 1.071114538622919
$ws.Range("E23").Value = 1.08015623413491
$ws.Range("F23").Value = 1.085554794325467
$ws.Range("I23").Value = 1.061028835606301
$ws.Range("J23").Value = 1.073552168343799
$ws.Range("K23").Value = 1.074627137182248
$ws.Range("L23").Value = 1.083637300935868
$ws.Range("M23").Value = 1.089017309823656
$ws.Range("N23").Value = 1.027957702259092

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.068737671472664
$ws.Range("D24").Value = 1.072444275793254
$ws.Range("E24").Value = 1.081738877887168
$ws.Range("F24").Value = 1.087111478311893
$ws.Range("I24").Value = 1.061625855375951
$ws.Range("J24").Value = 1.074714872275173
$ws.Range("K24").Value = 1.075691899025677
$ws.Range("L24").Value = 1.084956720631951
$ws.Range("M24").Value = 1.090312360842996
$ws.Range("N24").Value = 1.028365979865629

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.070671266246804
$ws.Range("D25").Value = 1.073987770275085
$ws.Range("E25").Value = 1.083577495431351
$ws.Range("F25").Value = 1.088919848744567
$ws.Range("I25").Value = 1.062313209250466
$ws.Range("J25").Value = 1.076061970734441
$ws.Range("K25").Value = 1.076924931986327
$ws.Range("L25").Value = 1.086487074409189
$ws.Range("M25").Value = 1.091814290324339
$ws.Range("N25").Value = 1.028837909218664
